$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("summary_stats")

# Update row 8 (Propriva) summary statistics
$ws.Range("B8").Value = 52.32092986
$ws.Range("C8").Value = 50.446266345
$ws.Range("D8").Value = 43.85201301302044
$ws.Range("E8").Value = 6.62208524658362
$ws.Range("F8").Value = 0.544332054194047
